$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Sheet1 -> Page1
$ws.Name = "Page1"

# Write the new cell content
$ws.Range("A1").Value = "1st change"

# Best-fit column A to its content, like double-clicking the column border
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the selection where the author left it
$ws.Range("F33").Select()
